$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Remove the now-obsolete per-VS-version "Auto-updates ... occur
#    successfully" rows for VS2005, VS2008, VS2010, VS2012 and VS2013
#    (original rows 76-80).
# ---------------------------------------------------------------------------
$ws.Rows("76:80").Delete()

# ---------------------------------------------------------------------------
# 2) The remaining "Auto-updates ... - VS2015" row (now row 76) is reworded
#    to cover "VS2015 & above" (it already was PASS).
# ---------------------------------------------------------------------------
$ws.Cells.Item(76, 2).Value = "Auto-updates (within IDE Extension Manager) occur successfully (VS2015 & above only)"
$ws.Cells.Item(76, 3).Value = "PASS"
$ws.Cells.Item(76, 3).Interior.Color = 5296274

# ---------------------------------------------------------------------------
# 3) "VS Gallery & AppVeyor version numbers in sync" (now row 77) flips from
#    TODO to PASS.
# ---------------------------------------------------------------------------
$ws.Cells.Item(77, 3).Value = "PASS"
$ws.Cells.Item(77, 3).Interior.Color = 5296274

# ---------------------------------------------------------------------------
# 4) Insert a new test-case row after row 39 ("When a rule exists for a
#    deeply nested file ...") for the new "directory not within the solution
#    file directory" test case. Inserting at row 40 pushes everything that
#    was there below it down by one.
# ---------------------------------------------------------------------------
$ws.Rows("40").Insert()

$ws.Cells.Item(40, 1).Value = "Functionality"
$ws.Cells.Item(40, 2).Value = "When a rule exists for a project in a directory that is not within the solution file directory, it is still modified successfully"
$ws.Cells.Item(40, 3).Value = "PASS"

# ---------------------------------------------------------------------------
# 5) Restore the selected cell in the bottom-right pane to B41 (matches the
#    saved view state in the edited workbook).
# ---------------------------------------------------------------------------
$ws.Range("B41").Select()
